$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the first three header cells (shared strings are re-created in this
# order: mass, wingarea, span) to match the target workbook.
$ws.Range("C1").Value = "mass"
$ws.Range("A1").Value = "wingarea"
$ws.Range("B1").Value = "span"

# Update the active selection to A2.
[void]$ws.Range("A2").Select()
